$wb = $excel.ActiveWorkbook

# --- Sheet "Weekly Quantity": append rows 14-16 ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("A13").Copy($ws1.Range("A14:A16"))
$ws1.Cells.Item(14, 1).Value = 45662.99999999999
$ws1.Cells.Item(14, 2).Value = 60
$ws1.Cells.Item(15, 1).Value = 45669.99999999999
$ws1.Cells.Item(15, 2).Value = 4
$ws1.Cells.Item(16, 1).Value = 45676.99999999999
$ws1.Cells.Item(16, 2).Value = 22

# --- Sheet "Monthly Trend": append row 8 ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("A7").Copy($ws2.Range("A8"))
$ws2.Cells.Item(8, 1).Value = 45688.99999999999
$ws2.Cells.Item(8, 2).Value = 86

# --- Sheet "PO Forecast": new forecast model, rows 2-24 ---
$ws3 = $wb.Worksheets.Item("PO Forecast")
$ws3.Range("A21").Copy($ws3.Range("A22:A24"))
$ws3.Cells.Item(2, 1).Value = 45494.99999999999
$ws3.Cells.Item(2, 2).Value = 55
$ws3.Cells.Item(3, 1).Value = 45522.99999999999
$ws3.Cells.Item(3, 2).Value = 56
$ws3.Cells.Item(4, 1).Value = 45557.99999999999
$ws3.Cells.Item(4, 2).Value = 57
$ws3.Cells.Item(5, 1).Value = 45564.99999999999
$ws3.Cells.Item(5, 2).Value = 57
$ws3.Cells.Item(6, 1).Value = 45571.99999999999
$ws3.Cells.Item(6, 2).Value = 57
$ws3.Cells.Item(7, 1).Value = 45578.99999999999
$ws3.Cells.Item(7, 2).Value = 58
$ws3.Cells.Item(8, 1).Value = 45585.99999999999
$ws3.Cells.Item(8, 2).Value = 58
$ws3.Cells.Item(9, 1).Value = 45592.99999999999
$ws3.Cells.Item(9, 2).Value = 58
$ws3.Cells.Item(10, 1).Value = 45599.99999999999
$ws3.Cells.Item(10, 2).Value = 58
$ws3.Cells.Item(11, 1).Value = 45606.99999999999
$ws3.Cells.Item(11, 2).Value = 58
$ws3.Cells.Item(12, 1).Value = 45613.99999999999
$ws3.Cells.Item(12, 2).Value = 58
$ws3.Cells.Item(13, 1).Value = 45634.99999999999
$ws3.Cells.Item(13, 2).Value = 59
$ws3.Cells.Item(14, 1).Value = 45662.99999999999
$ws3.Cells.Item(14, 2).Value = 60
$ws3.Cells.Item(15, 1).Value = 45669.99999999999
$ws3.Cells.Item(15, 2).Value = 60
$ws3.Cells.Item(16, 1).Value = 45676.99999999999
$ws3.Cells.Item(16, 2).Value = 60
$ws3.Cells.Item(17, 1).Value = 45683.99999999999
$ws3.Cells.Item(17, 2).Value = 60
$ws3.Cells.Item(18, 1).Value = 45690.99999999999
$ws3.Cells.Item(18, 2).Value = 60
$ws3.Cells.Item(19, 1).Value = 45697.99999999999
$ws3.Cells.Item(19, 2).Value = 61
$ws3.Cells.Item(20, 1).Value = 45704.99999999999
$ws3.Cells.Item(20, 2).Value = 61
$ws3.Cells.Item(21, 1).Value = 45711.99999999999
$ws3.Cells.Item(21, 2).Value = 61
$ws3.Cells.Item(22, 1).Value = 45718.99999999999
$ws3.Cells.Item(22, 2).Value = 61
$ws3.Cells.Item(23, 1).Value = 45725.99999999999
$ws3.Cells.Item(23, 2).Value = 61
$ws3.Cells.Item(24, 1).Value = 45732.99999999999
$ws3.Cells.Item(24, 2).Value = 62
